$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 1293
$ws.Range("F5").Value = 65
$ws.Range("F7").Value = 1002
$ws.Range("F11").Value = 109
$ws.Range("F12").Value = 437
$ws.Range("F14").Value = 1834
$ws.Range("F15").Value = 4252
$ws.Range("F18").Value = 2718
$ws.Range("F20").Value = 1106
$ws.Range("F21").Value = 3715
$ws.Range("F22").Value = 799
$ws.Range("F23").Value = 847
$ws.Range("F25").Value = 1511
$ws.Range("F27").Value = 122
$ws.Range("F28").Value = 873
$ws.Range("F29").Value = 179
$ws.Range("F31").Value = 239
$ws.Range("F33").Value = 31
$ws.Range("F34").Value = 1411
$ws.Range("F35").Value = 1996
$ws.Range("F36").Value = 948
$ws.Range("F37").Value = 9
$ws.Range("F38").Value = 516
$ws.Range("F39").Value = 88
$ws.Range("F41").Value = 602
$ws.Range("F42").Value = 300
$ws.Range("F43").Value = 118
$ws.Range("F45").Value = 244
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 154
$ws.Range("F9").Value = 19
$ws.Range("F12").Value = 124
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 481
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 481
$ws.Range("F4").Value = 1293
$ws.Range("F5").Value = 65
$ws.Range("F6").Value = 1002
$ws.Range("F8").Value = 154
$ws.Range("F13").Value = 109
$ws.Range("F15").Value = 437
$ws.Range("F16").Value = 1834
$ws.Range("F17").Value = 4252
$ws.Range("F21").Value = 2718
$ws.Range("F22").Value = 1106
$ws.Range("F23").Value = 3715
$ws.Range("F24").Value = 799
$ws.Range("F25").Value = 847
$ws.Range("F27").Value = 1511
$ws.Range("F29").Value = 19
$ws.Range("F32").Value = 122
$ws.Range("F33").Value = 124
$ws.Range("F34").Value = 873
$ws.Range("F35").Value = 179
$ws.Range("F37").Value = 239
$ws.Range("F39").Value = 1411
$ws.Range("F40").Value = 1996
$ws.Range("F41").Value = 948
$ws.Range("F42").Value = 516
$ws.Range("F43").Value = 88
$ws.Range("F44").Value = 602
$ws.Range("F45").Value = 300
$ws.Range("F46").Value = 118
$ws.Range("F48").Value = 244
